# This edit re-sorts the data rows (4-35) of the "Artfynd" sheet: each
# whole record (row) moves to a new row position, per the mapping below.
# No field-level values are changed - only row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns in the data table that hold text (inlineStr) values, including
# ones that look like numbers/dates/times (e.g. "1", "2014-09-05", "00:00").
# We force these to Text format BEFORE writing so Excel does not
# auto-convert them into real numbers/dates on assignment.
$textCols = @("C","D","F","G","H","I","J","K","L","M","N","P","T","U","V","W","Y","Z","AA","AB","AC","AF","AR","AT","AW","AX","AY")
foreach ($col in $textCols) {
    $ws.Range(($col + "4:" + $col + "35")).NumberFormat = "@"
}

# Snapshot the full current content (all columns A:AY) of every data row
# BEFORE any writes happen, so source data is never lost mid-shuffle.
$snaps = @{}
for ($r = 4; $r -le 35; $r++) {
    $snaps[$r] = $ws.Range(("A" + $r + ":AY" + $r)).Value2
}

# Mapping: destination row -> source row (content that should end up there)
$mapping = @{
    4 = 10
    5 = 34
    6 = 35
    7 = 4
    8 = 7
    9 = 9
    10 = 11
    11 = 22
    12 = 30
    13 = 31
    14 = 32
    15 = 33
    16 = 5
    17 = 6
    18 = 8
    19 = 12
    20 = 13
    21 = 14
    22 = 15
    23 = 16
    24 = 17
    25 = 18
    26 = 19
    27 = 20
    28 = 21
    29 = 23
    30 = 24
    31 = 25
    32 = 26
    33 = 27
    34 = 28
    35 = 29
}

foreach ($dest in $mapping.Keys) {
    $src = $mapping[$dest]
    $ws.Range(("A" + $dest + ":AY" + $dest)).Value2 = $snaps[$src]
}
